# Generate Report for Handoff
# The file 3d0ab4fa-ca32-4de2-af8c-b1e2b6657cd1.md is now "Ready for handoff":
#  - Overview sheet: zh-cn / de-de status columns + the "Latest HO Xliff
#    Generate Date" column for that row are refreshed.
#  - zh-cn / de-de sheets: Status, Priority and Latest Handoff Datetime for
#    that row are refreshed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet, row 3 (3d0ab4fa-ca32-4de2-af8c-b1e2b6657cd1.md) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 02:16:53"
$overview.Columns.Item(5).ColumnWidth = 17
$overview.Columns.Item(6).ColumnWidth = 17

# --- zh-cn sheet, row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-05 02:16:48"
$zhcn.Columns.Item(3).ColumnWidth = 17

# --- de-de sheet, row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-05 02:16:53"
$dede.Columns.Item(3).ColumnWidth = 17
